$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper pattern used throughout:
#   - Date-like text (e.g. "2022-09-16") is written through a temporary
#     formula ( ="2022-09-16" ) and then "Paste Values" (xlPasteValues = -4163)
#     over itself. This avoids Excel's automatic text->date conversion while
#     still landing on the plain General/style-5 formatting (a literal
#     Value="2022-09-16" assignment gets auto-parsed into a date serial).
#   - Numeric / plain-text cells are written directly via .Value.
#   - The "Run Time" (column B) cells reuse the existing datetime style from
#     a neighboring cell via PasteSpecial(xlPasteFormats = -4122).
# ---------------------------------------------------------------------------

# =====================================================================
# Sheet "AMSIN": append rows 45, 46, 47
# =====================================================================
$wsAMSIN = $wb.Worksheets.Item("AMSIN")

# --- Row 45 ---
$wsAMSIN.Cells.Item(45, 1).Formula = "=""2022-09-16"""
$wsAMSIN.Range("A45").Copy()
$wsAMSIN.Range("A45").PasteSpecial(-4163)
$wsAMSIN.Cells.Item(45, 2).Value = 44820.64129697916
$wsAMSIN.Cells.Item(45, 3).Value = "fstcyc167"
$wsAMSIN.Cells.Item(45, 4).Value = 119
$wsAMSIN.Cells.Item(45, 5).Value = 119
$wsAMSIN.Cells.Item(45, 6).Value = 0
$wsAMSIN.Cells.Item(45, 7).Value = 3.11
$wsAMSIN.Range("B44").Copy()
$wsAMSIN.Range("B45").PasteSpecial(-4122)

# --- Row 46 ---
$wsAMSIN.Cells.Item(46, 1).Formula = "=""2022-09-19"""
$wsAMSIN.Range("A46").Copy()
$wsAMSIN.Range("A46").PasteSpecial(-4163)
$wsAMSIN.Cells.Item(46, 2).Value = 44823.66725746528
$wsAMSIN.Cells.Item(46, 3).Value = "scndcycle167"
$wsAMSIN.Cells.Item(46, 4).Value = 119
$wsAMSIN.Cells.Item(46, 5).Value = 119
$wsAMSIN.Cells.Item(46, 6).Value = 0
$wsAMSIN.Cells.Item(46, 7).Value = 3.38
$wsAMSIN.Range("B44").Copy()
$wsAMSIN.Range("B46").PasteSpecial(-4122)

# --- Row 47 ---
$wsAMSIN.Cells.Item(47, 1).Formula = "=""2022-09-20"""
$wsAMSIN.Range("A47").Copy()
$wsAMSIN.Range("A47").PasteSpecial(-4163)
$wsAMSIN.Cells.Item(47, 2).Value = 44824.38552099537
$wsAMSIN.Cells.Item(47, 3).Value = "finalrun167"
$wsAMSIN.Cells.Item(47, 4).Value = 119
$wsAMSIN.Cells.Item(47, 5).Value = 118
$wsAMSIN.Cells.Item(47, 6).Value = 1
$wsAMSIN.Cells.Item(47, 7).Value = 3.09
$wsAMSIN.Range("B44").Copy()
$wsAMSIN.Range("B47").PasteSpecial(-4122)

# =====================================================================
# Sheet "BETA": append row 22 (unstyled cells except the Run Time column,
# matching the sparser styling already used on this sheet's latest rows)
# =====================================================================
$wsBETA = $wb.Worksheets.Item("BETA")

$wsBETA.Cells.Item(22, 1).Formula = "=""2022-09-20"""
$wsBETA.Range("A22").Copy()
$wsBETA.Range("A22").PasteSpecial(-4163)

$wsBETA.Cells.Item(22, 2).Value = 44824.53252232823
$wsBETA.Range("B21").Copy()
$wsBETA.Range("B22").PasteSpecial(-4122)

$wsBETA.Cells.Item(22, 3).Value = "beta167"
$wsBETA.Cells.Item(22, 4).Value = 119
$wsBETA.Cells.Item(22, 5).Value = 119
$wsBETA.Cells.Item(22, 6).Value = 0
$wsBETA.Cells.Item(22, 7).Value = 3.27

# =====================================================================
# Sheet "AMS": row 22 gains the s="5" formatting used elsewhere on the
# sheet plus a corrected Run Time (B22) value.
# =====================================================================
$wsAMS = $wb.Worksheets.Item("AMS")

$wsAMS.Range("A22:G22").ClearContents()

$wsAMS.Cells.Item(22, 1).Formula = "=""2022-08-24"""
$wsAMS.Range("A22").Copy()
$wsAMS.Range("A22").PasteSpecial(-4163)

$wsAMS.Cells.Item(22, 2).Value = 44797.92940078703

$wsAMS.Cells.Item(22, 3).Value = "166_live"
$wsAMS.Cells.Item(22, 4).Value = 119
$wsAMS.Cells.Item(22, 5).Value = 119
$wsAMS.Cells.Item(22, 6).Value = 0
$wsAMS.Cells.Item(22, 7).Value = 2.8
